$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.774.21'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +7.49%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.596.26'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.33%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '418.32'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.48%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.86'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.53%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.651'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.00%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.589.59'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.44%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.13%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.770'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.81%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +16.63%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000343'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +51.50%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.64'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.50%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.99'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.82%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.144.35'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.40%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.23%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.49'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.82%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.590.09'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.21%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +4.71%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.597.55'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +7.26%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.39'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.56%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '463.31'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.11%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '88.64'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.50%  '

# Row 24
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.35%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.44'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.05%  '

# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Filecoin'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.27'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.98%  '

# Row 27
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.37'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.60%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '35.55'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.71%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.86'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.35%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.79'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.21%  '

# Row 31
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.42%  '

# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.46'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.71%  '

# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.118'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.34%  '

# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.163'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.99%  '

# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.67'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.23%  '

# Row 36
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.10%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.85'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.90%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0496'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.80%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0718'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +23.05%  '

# Row 40
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +7.91%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.23%  '

# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.01%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '148.46'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.16%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.73'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.53%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -2.29%  '

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.26%  '

# Row 47
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.12%  '

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.37%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.67%  '

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +16.57%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '15.75'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.25%  '
